# Weekly fruit/vegetable price update: a new week's price record is
# inserted at the top of the data history (row 36), pushing all the
# existing historical rows (36-60) down by one (to 37-61).
#
# The new row 36 carries the same static descriptive fields (mercado,
# region, category, quality, unit, origin, classification, etc.) as the
# series already uses, with a new date and new min/max/avg/kg prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 36; everything below shifts down one.
$ws.Rows.Item(36).Insert()

# Populate the newly inserted row 36 with this week's record.
$ws.Range("A36").Value = 1
$ws.Range("B36").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C36").Value = "Arica y Parinacota"
$ws.Range("D36").Value = 44566
$ws.Range("E36").Value = 15
$ws.Range("F36").Value = 100112038
$ws.Range("G36").Value = "Cebollín baby"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 250
$ws.Range("K36").Value = 4000
$ws.Range("L36").Value = 4500
$ws.Range("M36").Value = 4250
$ws.Range("N36").Value = "$/paquete 1,5 a 2 kilos"
$ws.Range("O36").Value = "Región de Arica y Parinacota"
$ws.Range("P36").Value = 2125
$ws.Range("Q36").Value = 2
$ws.Range("R36").Value = "Hortaliza"
